# feat: add 2022-Q3 data
#
# The workbook currently has:
#   Sheet 1: "总计"      (summary)
#   Sheet 2: "2022-Q2"   (quarterly fund-holding detail)
#
# We need to:
#   1. Insert a brand-new "2022-Q3" quarter as the summary's 2nd data row
#      (pushing the existing "2022-Q2" summary row down to row 3).
#   2. Turn the existing "2022-Q2" detail sheet into the new "2022-Q3"
#      detail sheet (new numbers), while keeping a duplicate of the
#      original "2022-Q2" data (old numbers) as its own sheet placed
#      right after it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "总计" summary sheet: shift the existing data row down and add
#    the new 2022-Q3 row above it.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Duplicate row 2's formatting into row 3 (keeps style s="2" on A3,
# matching the existing A2), then fill in the old 2022-Q2 values that
# used to live in row 2.
$summary.Range("A2").Copy($summary.Range("A3"))
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.02

# Row 2 becomes the new 2022-Q3 entry.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("D2").Value = 0.01

# ---------------------------------------------------------------
# 2) Detail sheets: duplicate the current "2022-Q2" sheet so the
#    original data is preserved on its own tab, then overwrite the
#    original sheet in place with the new 2022-Q3 figures.
# ---------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)

# Make an exact copy (values + styles) right after the original sheet;
# this copy keeps the old data (Excel auto-names it "2022-Q2 (2)" since
# the source name is still taken at copy time).
$q2.Copy($null, $q2)
$copy = $wb.Worksheets.Item(3)

# Free up the "2022-Q2" name, then rename the original sheet to
# 2022-Q3 and the copy back to 2022-Q2.
$q2.Name = "2022-Q3"
$copy.Name = "2022-Q2"

# The new 2022-Q3 sheet uses the same header/label styling as the
# "总计" sheet (style s="2") rather than the old per-quarter sheet's
# own header style (s="1") that it inherited from the rename above;
# copy just the formatting over, leaving all values untouched.
$summary.Range("B1").Copy()
$q2.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q2.Range("A2:A3").PasteSpecial(-4122)

# D:G hold numeric-looking figures that are stored as plain TEXT (not
# numbers) in this workbook's convention, same as the existing sheet.
# A leading "'" forces literal text entry instead of Excel's automatic
# number parsing; the style is reset right after so the quote-prefix
# formatting doesn't stick around on the cells.
$q2.Range("D2").Value = "'0.67"
$q2.Range("E2").Value = "'90.27"
$q2.Range("F2").Value = "'1.62"
$q2.Range("G2").Value = "'0.0109"
$q2.Range("H2").Value = 7

$q2.Range("D3").Value = "'0.06"
$q2.Range("E3").Value = "'90.27"
$q2.Range("F3").Value = "'1.62"
$q2.Range("G3").Value = "'0.0010"
$q2.Range("H3").Value = 7

$q2.Range("D2:G3").Style = "Normal"
